$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 previously held "xls with qa.guru"; replace with the new text.
# B5 already reads "some text" and stays that way (Excel's shared-string
# table re-orders itself as a side effect of this edit).
$ws.Range("B4").Value = "hello qa.guru students!"
$ws.Range("B5").Value = "some text"

# Recolor the JetBrains Mono font used by B4 from the old green to the new gray.
$ws.Range("B4").Font.Color = 13023145

# Update the sheet's remembered selection/active cell.
$ws.Range("F7").Select()
